# Refresh the crypto price/volume snapshot (cryptos list update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.536.63"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.685.05"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").Value = "'313.89"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "'0.3897"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("D8").Value = "'0.4016"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").Value = "'1.479"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "'52.91"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'0.08682"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").Value = "'7.627"
$ws.Range("E13").Value = "  +5.76%  "
$ws.Range("D14").Value = "'24.40"
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("D15").Value = "'7.925"
$ws.Range("E15").Value = "  -1.42%  "
$ws.Range("D16").Value = "'0.00001328"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "1.684.25"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'97.89"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Value = "'0.07090"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").Value = "'19.57"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "'7.243"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").Value = "24.524.43"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "'2.985"
$ws.Range("E25").Value = "  -8.97%  "
$ws.Range("D26").Value = "'2.349"
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "'22.55"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "'161.13"
$ws.Range("E28").Value = "  -0.91%  "
$ws.Range("D29").Value = "'8.465"
$ws.Range("E29").Value = "  +11.29%  "
$ws.Range("D30").Value = "'5.239"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("D31").Value = "'136.17"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "1.870.21"
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "'7.443"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").Value = "'0.08715"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "'1.027"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").Value = "'1.948"
$ws.Range("E36").Value = "  +2.74%  "
$ws.Range("D37").Value = "'0.02876"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("D38").Value = "'0.2698"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "'10.64"
$ws.Range("E39").Value = "  -4.50%  "
$ws.Range("D40").Value = "'0.09085"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'13.98"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "'0.7690"
$ws.Range("E42").Value = "  +1.42%  "
$ws.Range("D43").Value = "'1.447"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "'16.50"
$ws.Range("E44").Value = "  +2.90%  "
$ws.Range("D45").Value = "'0.7088"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "'2.550"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "'4.195"
$ws.Range("E47").Value = "  -0.45%  "
$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'1.327"
$ws.Range("D50").Value = "'137.81"
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "'90.52"
$ws.Range("E51").Value = "  +1.28%  "
